$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.855.31'
$ws.Range('E2').Value = '  +1.28%  '
$ws.Range('D3').Value = '1.812.12'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').Value = '225.12'
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('D6').Value = '0.606'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('D8').Value = '39.49'
$ws.Range('E8').Value = '  +9.08%  '
$ws.Range('D9').Value = '0.291'
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('E10').Value = '  -2.92%  '
$ws.Range('D11').Value = '0.0997'
$ws.Range('E11').Value = '  +3.22%  '
$ws.Range('D12').Value = '2.075.40'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').Value = '1.812.87'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').Value = '11.03'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').Value = '34.867.34'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('D17').Value = '4.41'
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').Value = '68.29'
$ws.Range('E18').Value = '  -2.55%  '
$ws.Range('D19').Value = '241.79'
$ws.Range('E19').Value = '  -1.49%  '
$ws.Range('D20').Value = '0.0₃0771'
$ws.Range('E20').Value = '  -2.14%  '
$ws.Range('D21').Value = '11.18'
$ws.Range('E21').Value = '  -2.64%  '
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').Value = '4.11'
$ws.Range('E23').Value = '  -1.70%  '
$ws.Range('D24').Value = '2.17'
$ws.Range('E24').Value = '  -4.10%  '
$ws.Range('D25').Value = '171.45'
$ws.Range('D26').Value = '7.76'
$ws.Range('E26').Value = '  -4.31%  '
$ws.Range('D27').Value = '17.56'
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('E28').Value = '  -0.45%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '3.78'
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.22'
$ws.Range('E31').Value = '  -1.46%  '
$ws.Range('D32').Value = '3.88'
$ws.Range('E32').Value = '  -2.96%  '
$ws.Range('D33').Value = '0.0516'
$ws.Range('E33').Value = '  -1.50%  '
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '0.647'
$ws.Range('E35').Value = '  -2.23%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '1.316.18'
$ws.Range('E36').Value = '  -4.79%  '
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('E38').Value = '  +1.10%  '
$ws.Range('D39').Value = '0.0189'
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('D40').Value = '1.25'
$ws.Range('E40').Value = '  +5.48%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = '14.70'
$ws.Range('E41').Value = '  +8.67%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '83.10'
$ws.Range('E42').Value = '  +0.88%  '
$ws.Range('D43').Value = '2.43'
$ws.Range('E43').Value = '  +0.46%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '2.81'
$ws.Range('E44').Value = '  -0.29%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').Value = '0.951'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').Value = '0.0518'
$ws.Range('E46').Value = '  +3.70%  '
$ws.Range('D47').Value = '1.977.42'
$ws.Range('D48').Value = '5.77'
$ws.Range('E48').Value = '  -3.36%  '
$ws.Range('D49').Value = '0.997'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('D50').Value = '102.42'
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '48.64'
$ws.Range('E51').Value = '  -1.32%  '
